$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "param_TimeStep_starting_index",
    "param_demand1_op_cost_starting_index",
    "param_demand1_inv_cost_starting_index",
    "param_net1_buy_electric_starting_index",
    "param_P_from_net1_starting_index",
    "param_net1_sell_electric_starting_index",
    "param_Q_from_net1_starting_index",
    "param_net1_buy_thermal_starting_index",
    "param_P_net1_demand1_starting_index",
    "param_net1_sell_thermal_starting_index",
    "param_Q_to_net1_starting_index",
    "param_P_to_net1_starting_index",
    "param_net1_emissions_starting_index",
    "param_net1_inv_cost_starting_index",
    "param_Q_net1_demand1_starting_index",
    "param_P_net1_bat1_starting_index",
    "param_pv1_op_cost_starting_index",
    "param_P_from_pv1_starting_index",
    "param_P_pv1_net1_starting_index",
    "param_pv1_inv_cost_starting_index",
    "param_pv1_emissions_starting_index",
    "param_P_pv1_demand1_starting_index",
    "param_P_pv1_bat1_starting_index",
    "param_P_from_bat1_starting_index",
    "param_bat1_emissions_starting_index",
    "param_bat1_cumulated_aging_starting_index",
    "param_bat1_K_dis_starting_index",
    "param_P_bat1_net1_starting_index",
    "param_bat1_K_ch_starting_index",
    "param_bat1_SOC_max_starting_index",
    "param_bat1_inv_cost_starting_index",
    "param_P_bat1_demand1_starting_index",
    "param_bat1_SOC_starting_index",
    "param_bat1_op_cost_starting_index",
    "param_bat1_integer_starting_index",
    "param_P_to_bat1_starting_index",
    "param_total_operation_cost_starting_index",
    "param_total_buy_starting_index",
    "param_total_sell_starting_index",
    "param_total_emissions_starting_index"
)

$values = @(
    20,
    0,
    0,
    177.7435074899564,
    433.5207499755035,
    0,
    650.2811249632553,
    214.5927712378742,
    433.5207499755035,
    0,
    0,
    0,
    398.8390899774632,
    0,
    650.2811249632553,
    0,
    1,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0,
    0.0001088948573468242,
    1,
    0,
    0,
    0.9998911051426532,
    0,
    0,
    0.9998911051426532,
    1,
    -0,
    0,
    2,
    392.3362787278307,
    0,
    398.8390899774632
)

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}